$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing subscription row (was "MCAPS-MarcusGaspar") to the new subscription name
$ws.Range("A2").Value = "Default - Microsoft Azure Sponsorship 2"
$ws.Range("B2").Value = "FinOps3"

# Add a new row for the second subscription that also gets the FinOps3 tag
$ws.Range("A3").Value = "HPC subscription"
$ws.Range("B3").Value = "FinOps3"

# Column A needs to widen to fit the longer subscription names (re-run best-fit)
$ws.Columns.Item(1).ColumnWidth = 33

# Move the active selection to A6, matching the saved view state
$ws.Range("A6").Select() | Out-Null
